# Fix the typo in the shared string for cell A7 ("Pn_ subpacifica_B" -> "Pn_subpacifica_B")
# and update the active selection to A7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Pn_subpacifica_B"
$ws.Range("A7").Select()
